$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.147712230682373
$ws.Range("B1").Value = 7.166317462921143
$ws.Range("C1").Value = 5.629353523254395
$ws.Range("D1").Value = -1
$ws.Range("E1").Value = 3.168505191802979
